$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.451.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.108.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("E7").Value = "  +1.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08917"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.090.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.798"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.006"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.93%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001132"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06629"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.334"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.509.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.366"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.346.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.583"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.225"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.97%  "

$ws.Range("E32").Value = "  +15.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1074"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.193"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.926"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.06%  "

$ws.Range("E37").Value = "  +1.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06842"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.559"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2308"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6915"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.84%  "

$ws.Range("E43").Value = "  +1.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.353"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6384"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.666"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000355"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +26.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.248"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.03%  "
